# Update Thresholds and Results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column width adjustments ---
# Column B widens slightly (existing bestFit column); columns D, E, G, H
# pick up explicit widths that previously relied on the default width.
$ws.Columns.Item(2).ColumnWidth = 17.25
$ws.Columns.Item(4).ColumnWidth = 14.45
$ws.Columns.Item(5).ColumnWidth = 11.25
$ws.Columns.Item(7).ColumnWidth = 14.6
$ws.Columns.Item(8).ColumnWidth = 11.25

# --- Updated threshold / result values ---
$ws.Range("E2").Value = 0.78225367446924332
$ws.Range("F2").Value = 0.94
$ws.Range("G2").Value = 0.89112683723462172
$ws.Range("H2").Value = 0.81818181818181823

$ws.Range("D3").Value = 0.92999999999999994
$ws.Range("E3").Value = 0.85689896326131898
$ws.Range("G3").Value = 0.89714285714285724
$ws.Range("H3").Value = 0.86956521739130432

$ws.Range("E4").Value = 0.80995040343199776
$ws.Range("F4").Value = 0.96499999999999997
$ws.Range("H4").Value = 0.81081081081081086

$ws.Range("D5").Value = 0.91999999999999993
$ws.Range("E5").Value = 0.97974678822934125
$ws.Range("F5").Value = 0.995
$ws.Range("G5").Value = 0.99709302325581395
$ws.Range("H5").Value = 0.98245614035087714

$ws.Range("D6").Value = 0.73
$ws.Range("E6").Value = 0.9691664180319719
$ws.Range("F6").Value = 0.995
$ws.Range("G6").Value = 0.99726775956284153
$ws.Range("H6").Value = 0.97142857142857142

$ws.Range("D7").Value = 0.92999999999999994
$ws.Range("E7").Value = 0.96349126212248548
$ws.Range("F7").Value = 0.995
$ws.Range("G7").Value = 0.99731182795698925
$ws.Range("H7").Value = 0.96551724137931039

$ws.Range("E8").Value = 0.50009913796902061
$ws.Range("F8").Value = 0.85499999999999998
$ws.Range("G8").Value = 0.77745098039215677
$ws.Range("H8").Value = 0.57971014492753625

$ws.Range("D9").Value = 0.86999999999999988
$ws.Range("E9").Value = 0.79633669583420641
$ws.Range("F9").Value = 0.94
$ws.Range("G9").Value = 0.95016611295681064
$ws.Range("H9").Value = 0.81818181818181823

$ws.Range("E11").Value = 0.93298784624641462
$ws.Range("F11").Value = 0.995
$ws.Range("G11").Value = 0.99740932642487046
$ws.Range("H11").Value = 0.93333333333333335
